$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (D: price, E: percentage) are stored as literal text,
# matching the source data which encodes these columns as plain strings.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '278.81'
$ws.Cells.Item(2, 5).Value = '6.79%'

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '27.20'
$ws.Cells.Item(3, 5).Value = '0.68%'

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '4.821'
$ws.Cells.Item(4, 5).Value = '2.79%'

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.06264'
$ws.Cells.Item(5, 5).Value = '0.74%'

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '6.862'
$ws.Cells.Item(6, 5).Value = '1.70%'

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = 'MXToken'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(7, 4).Value = '0.8788'
$ws.Cells.Item(7, 5).Value = '3.15%'

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = 'FTXToken'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(8, 4).Value = '0.9390'
$ws.Cells.Item(8, 5).Value = '2.74%'

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = 'WazirX'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(9, 4).Value = '0.1449'
$ws.Cells.Item(9, 5).Value = '3.52%'

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(10, 4).Value = '0.05155'
$ws.Cells.Item(10, 5).Value = '6.27%'

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(11, 4).Value = '0.07284'
$ws.Cells.Item(11, 5).Value = '2.77%'

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = 'BitrueCoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(12, 4).Value = '0.03162'
$ws.Cells.Item(12, 5).Value = '1.77%'

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = 'BitMartToken'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(13, 4).Value = '0.09054'
$ws.Cells.Item(13, 5).Value = '-0.09%'

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = 'BitForexToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(14, 4).Value = '0.001562'
$ws.Cells.Item(14, 5).Value = '1.33%'

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = 'One'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(15, 4).Value = '0.0006277'
$ws.Cells.Item(15, 5).Value = '1.48%'

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = 'TigerCash'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(16, 4).Value = '0.005984'
$ws.Cells.Item(16, 5).Value = '-0.49%'

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Cells.Item(17, 2).Value = 'LEO'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(17, 4).Value = '3.450'
$ws.Cells.Item(17, 5).Value = '0.25%'

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Cells.Item(18, 2).Value = 'GateToken'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(18, 4).Value = '3.265'
$ws.Cells.Item(18, 5).Value = '2.81%'

$ws.Range("E19").NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '5.56%'

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.1310'
$ws.Cells.Item(21, 5).Value = '-0.05%'

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '3.849'
$ws.Cells.Item(22, 5).Value = '-5.92%'

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.04317'
$ws.Cells.Item(23, 5).Value = '1.44%'

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.001175'
$ws.Cells.Item(24, 5).Value = '-2.74%'

$ws.Range("E25").NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '4.81%'

$ws.Range("E26").NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '-0.19%'

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.0001689'
$ws.Cells.Item(27, 5).Value = '2.99%'

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.04036'
$ws.Cells.Item(40, 5).Value = '2.32%'

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.006399'
$ws.Cells.Item(41, 5).Value = '55.52%'

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.1153'
$ws.Cells.Item(42, 5).Value = '3.65%'

$ws.Range("E43").NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '-4.89%'

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.01394'
$ws.Cells.Item(44, 5).Value = '0.44%'

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.00005152'
$ws.Cells.Item(45, 5).Value = '-0.21%'

$ws.Range("E46").NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '-0.14%'

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.362'
$ws.Cells.Item(47, 5).Value = '1,023.35%'

$ws.Range("E49").NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '-0.14%'

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0001999'
$ws.Cells.Item(50, 5).Value = '-0.14%'

